# DYK banner slide: swap the "#4 Learn more..." call-to-action subtitle
# for the new "#9 Roadmap Insights" series teaser (new DYK series videos).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the subtitle placeholder shape (named "Subtitle 15", id 16) that
# currently reads "#4 Learn more..."
$sh = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Name -eq "Subtitle 15") {
        $sh = $cand
    }
}
if ($sh -eq $null) {
    $sh = $s.Shapes.Item(2)
}

$tr = $sh.TextFrame.TextRange

# Build the full replacement text as a single run first.
$tr.Text = "#9 Roadmap Insights"

# The second run ("Roadmap Insights") keeps the paragraph's original
# en-GB language; stamp that onto the whole range before splitting so the
# split inherits it cleanly on both sides.
$tr.LanguageID = "en-GB"

# Split "#9 " into its own run (a same-value Text assignment on the
# sub-range forces PowerPoint to break it out as a separate run without
# touching any other formatting).
$lead = $tr.Characters(1, 3)
$lead.Text = "#9 "

# Now restore the first run ("#9 ") back to the en-US language the title
# run originally used; this targets the paragraph's first run, i.e.
# exactly the "#9 " run we just split off, leaving "Roadmap Insights" as
# en-GB.
$tr.LanguageID = "en-US"
